$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.951.56"
$ws.Range("E2").Value = "  +8.06%  "

$ws.Range("D3").Value = "3.634.11"
$ws.Range("E3").Value = "  +4.49%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'418.71"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("D6").Value = "'133.10"
$ws.Range("E6").Value = "  +2.40%  "

$ws.Range("D7").Value = "'0.648"
$ws.Range("E7").Value = "  +3.50%  "

$ws.Range("D8").Value = "3.629.34"
$ws.Range("E8").Value = "  +4.65%  "

$ws.Range("D9").Value = "'0.998"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").Value = "'0.770"
$ws.Range("E10").Value = "  +5.96%  "

$ws.Range("D11").Value = "'0.182"
$ws.Range("E11").Value = "  +18.30%  "

$ws.Range("D12").Value = "'0.0000352"
$ws.Range("E12").Value = "  +56.55%  "

$ws.Range("D13").Value = "'42.83"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").Value = "'9.90"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").Value = "4.190.82"
$ws.Range("E15").Value = "  +3.88%  "

$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "'20.38"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").Value = "3.629.57"
$ws.Range("E18").Value = "  +4.37%  "

$ws.Range("D19").Value = "'1.15"
$ws.Range("E19").Value = "  +5.35%  "

$ws.Range("D20").Value = "67.807.04"
$ws.Range("E20").Value = "  +7.85%  "

$ws.Range("D21").Value = "'12.38"
$ws.Range("E21").Value = "  -1.88%  "

$ws.Range("D22").Value = "'463.42"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("D23").Value = "'88.59"
$ws.Range("E23").Value = "  -2.19%  "

$ws.Range("D24").Value = "'3.13"
$ws.Range("E24").Value = "  -5.10%  "

$ws.Range("D25").Value = "'13.37"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("D26").Value = "'3.35"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'35.74"
$ws.Range("E27").Value = "  +7.04%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = "  -4.36%  "

$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("E30").Value = "  +3.94%  "

$ws.Range("D31").Value = "'12.37"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "'7.41"
$ws.Range("E32").Value = "  -2.00%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.117"
$ws.Range("E33").Value = "  +4.03%  "

$ws.Range("E34").Value = "  -3.44%  "

$ws.Range("D35").Value = "'40.87"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "'56.75"
$ws.Range("E37").Value = "  -2.13%  "

$ws.Range("E38").Value = "  +1.12%  "

$ws.Range("D39").Value = "0.0₃0704"
$ws.Range("E39").Value = "  +20.59%  "

$ws.Range("D40").Value = "'0.146"
$ws.Range("E40").Value = "  +7.82%  "

$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.31%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "'148.05"
$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("E44").Value = "  -4.15%  "

$ws.Range("D45").Value = "'3.28"
$ws.Range("E45").Value = "  -1.19%  "

$ws.Range("D46").Value = "'4.32"
$ws.Range("E46").Value = "  -2.23%  "

$ws.Range("D47").Value = "'0.309"
$ws.Range("E47").Value = "  -3.49%  "

$ws.Range("D48").Value = "'1.98"
$ws.Range("E48").Value = "  -3.00%  "

$ws.Range("D49").Value = "'2.34"
$ws.Range("E49").Value = "  -1.79%  "

$ws.Range("D50").Value = "'2.70"
$ws.Range("E50").Value = "  +15.82%  "

$ws.Range("D51").Value = "'115.22"
$ws.Range("E51").Value = "  +5.88%  "
